$d = $word.ActiveDocument

function Find-Range($searchText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $rng.Find.Found) {
        throw "Not found: $searchText"
    }
    return $rng
}

function Color-Text($searchText, $color) {
    $rng = Find-Range $searchText
    $rng.Font.Color = $color
}

function Color-Prefix($searchText, $prefixLen, $color) {
    $rng = Find-Range $searchText
    $sub = $d.Range($rng.Start, $rng.Start + $prefixLen)
    $sub.Font.Color = $color
}

function Insert-After($searchText, $insertText, $color) {
    $rng = Find-Range $searchText
    $rng.Collapse(0)
    $ins = $d.Range($rng.Start, $rng.Start)
    $ins.InsertAfter($insertText)
    if ($color -ne $null) {
        $ins.Font.Color = $color
    }
}

$wdRed = 255

# 1. "0-4-9-1" -> recolor red, then append " (xoang)" red
Color-Text "0-4-9-1" $wdRed
Insert-After "0-4-9-1" " (xoang)" $wdRed

# 2. ", 0-6-4-1" -> recolor red, then append " (rìa)" red
Color-Text ", 0-6-4-1" $wdRed
Insert-After "0-6-4-1" " (rìa)" $wdRed

# 3. ", 0-10-7-1" -> recolor red, then append "(xoang)" red
Color-Text ", 0-10-7-1" $wdRed
Insert-After "0-10-7-1" "(xoang)" $wdRed

# 4. ", 1-7-11-1" stays black; append "(bỏ qua)" black (no color change)
Insert-After "1-7-11-1" "(bỏ qua)" $null

# 5. "1-9-14-1" (token only, not the leading ", ") -> recolor red, append " (xoang)" red
Color-Text "1-9-14-1" $wdRed
Insert-After "1-9-14-1" " (xoang)" $wdRed

# 6. The ", " right after "(xoang)" (before "1-10-4-0") -> recolor red (only the ", " prefix)
Color-Prefix ", 1-10-4-0" 2 $wdRed

# 7. Append " (dư)" black after "1-10-4-0" (no color change)
Insert-After "1-10-4-0" " (dư)" $null

# 8. ", 1-14-18-1" stays black; append " (thua" black, then ")" red
Insert-After "1-14-18-1" " (thua" $null
Insert-After " (thua" ")" $wdRed

# 9. ", 1-25-14-1" -> recolor red; append " (sáng quá)" red
Color-Text ", 1-25-14-1" $wdRed
Insert-After "1-25-14-1" " (sáng quá)" $wdRed

# 10. ", 1-30-17-1" -> recolor red; append " (sáng quá)" red
Color-Text ", 1-30-17-1" $wdRed
Insert-After "1-30-17-1" " (sáng quá)" $wdRed

# 11. ", 1-31-14-1" stays black; append "(bỏ qua)" black (no color change)
Insert-After "1-31-14-1" "(bỏ qua)" $null

# 12. Reposition the "_GoBack" bookmark to wrap just "1-9-14-1 (xoang), "
#     (mirrors Word's real behaviour of parking _GoBack at the most recent edit).
$bmRng = Find-Range "1-9-14-1 (xoang), "
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null
